$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (the "Subject" column) repeats the same text down every data row.
# Fix the typo/trailing-space: "Biochemistry Lab/CBL " -> "Biochemistry LAB/CBL".
$xlUp = -4162
$lastRow = $ws.Cells($ws.Rows.Count, 2).End($xlUp).Row
$ws.Range($ws.Cells(2, 2), $ws.Cells($lastRow, 2)).Value = "Biochemistry LAB/CBL"

# Move the active selection from the header row (A1:XFD1) to B2.
$ws.Range("B2").Select()
